$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 94.8
$ws.Range("J38").Value = 500
$ws.Range("L38").Value = 1500
$ws.Range("N38").Value = -2244

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 9543.286
$ws.Range("J37").Value = 10717.167
$ws.Range("L37").Value = 10717.167
$ws.Range("N37").Value = -11263.167
$ws.Range("H80").Value = 17150
$ws.Range("J80").Value = 28500
$ws.Range("L80").Value = 28500
$ws.Range("N80").Value = -30496
$ws.Range("H83").Value = 17150
$ws.Range("J83").Value = 28500
$ws.Range("L83").Value = 85500
$ws.Range("N83").Value = -95484
$ws.Range("H102").Value = 1825.579
$ws.Range("I102").Value = 1652
$ws.Range("K102").Value = 1652
$ws.Range("M102").Value = -30
$ws.Range("H121").Value = 25572.857
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 25572.857
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 25572.857
$ws.Range("N121").Value = -29066.857
$ws.Range("H122").Value = 1426.931
$ws.Range("I122").Value = 1415.7727
$ws.Range("J122").Value = 1462
$ws.Range("K122").Value = 4247.3181
$ws.Range("L122").Value = 4386
$ws.Range("M122").Value = -1797.3181
$ws.Range("N122").Value = -9286
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H124").Value = 12999.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 12999.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 12999.5
$ws.Range("N124").Value = -22819.5
$ws.Range("H125").Value = 42000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 42000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 42000
$ws.Range("N125").Value = -51840
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H127").Value = 52979.668
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 52979.668
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 52979.668
$ws.Range("N127").Value = -62899.668
$ws.Range("H128").Value = 36214.5
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 36214.5
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 36214.5
$ws.Range("N128").Value = -46174.5
$ws.Range("H129").Value = 37999.6
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 37999.6
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 37999.6
$ws.Range("N129").Value = -47999.6
$ws.Range("H130").Value = 27500
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 27500
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 27500
$ws.Range("N130").Value = -37540
$ws.Range("H131").Value = 46815
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 46815
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 46815
$ws.Range("N131").Value = -56895
$ws.Range("H132").Value = 37087.9
$ws.Range("I132").Value = 2763.5
$ws.Range("J132").Value = 93255.09
$ws.Range("K132").Value = 8290.5
$ws.Range("L132").Value = 279765.27
$ws.Range("M132").Value = -5760.5
$ws.Range("N132").Value = -284825.27
$ws.Range("H133").Value = 15000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 15000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 15000
$ws.Range("N133").Value = -20060
$ws.Range("H134").Value = 40000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 40000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 40000
$ws.Range("N134").Value = -50140
$ws.Range("H135").Value = 40518.43
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 40518.43
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 40518.43
$ws.Range("N135").Value = -50658.43
$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 50000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
$ws.Range("H138").Value = 38299.715
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 38299.715
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 38299.715
$ws.Range("N138").Value = -48579.715
$ws.Range("H139").Value = 35560.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 35560.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 35560.75
$ws.Range("N139").Value = -45840.75
$ws.Range("H140").Value = 37849.168
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 37849.168
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 37849.168
$ws.Range("N140").Value = -48209.168
$ws.Range("H141").Value = 41196.875
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 41196.875
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 41196.875
$ws.Range("N141").Value = -51556.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14172.429
$ws.Range("I82").Value = 5973.4287
$ws.Range("J82").Value = 22371.428
$ws.Range("K82").Value = 5973.4287
$ws.Range("L82").Value = 22371.428
$ws.Range("M82").Value = -5590.4287
$ws.Range("N82").Value = -23137.428
$ws.Range("H85").Value = 14172.429
$ws.Range("I85").Value = 5973.4287
$ws.Range("J85").Value = 22371.428
$ws.Range("K85").Value = 5973.4287
$ws.Range("L85").Value = 22371.428
$ws.Range("M85").Value = -4647.4287
$ws.Range("N85").Value = -25023.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9883.9
$ws.Range("J50").Value = 10353.714
$ws.Range("L50").Value = 10353.714
$ws.Range("N50").Value = -11603.714
$ws.Range("H51").Value = 10229.3
$ws.Range("J51").Value = 11022.556
$ws.Range("L51").Value = 11022.556
$ws.Range("N51").Value = -12494.556
$ws.Range("H61").Value = 10229.3
$ws.Range("J61").Value = 11022.556
$ws.Range("L61").Value = 11022.556
$ws.Range("N61").Value = -11718.556
$ws.Range("H68").Value = 18195.428
$ws.Range("J68").Value = 20300
$ws.Range("L68").Value = 20300
$ws.Range("N68").Value = -21798
$ws.Range("H71").Value = 18195.428
$ws.Range("J71").Value = 20300
$ws.Range("L71").Value = 60900
$ws.Range("N71").Value = -68388
$ws.Range("H74").Value = 16618.572
$ws.Range("J74").Value = 16618.572
$ws.Range("L74").Value = 16618.572
$ws.Range("N74").Value = -18366.572
$ws.Range("H77").Value = 16618.572
$ws.Range("J77").Value = 16618.572
$ws.Range("L77").Value = 49855.716
$ws.Range("N77").Value = -58591.716
$ws.Range("H81").Value = 42295.2
$ws.Range("I81").Value = 15000
$ws.Range("K81").Value = 15000
$ws.Range("M81").Value = -14002
$ws.Range("H84").Value = 42295.2
$ws.Range("I84").Value = 15000
$ws.Range("K84").Value = 45000
$ws.Range("M84").Value = -40008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 15380.857
$ws.Range("I120").Value = 8000
$ws.Range("J120").Value = 15948.615
$ws.Range("K120").Value = 24000
$ws.Range("L120").Value = 47845.845
$ws.Range("M120").Value = -19162
$ws.Range("N120").Value = -57521.845
$ws.Range("H121").Value = 35715140
$ws.Range("I121").Value = 530
$ws.Range("J121").Value = 62501100
$ws.Range("K121").Value = 1590
$ws.Range("L121").Value = 187503300
$ws.Range("M121").Value = -280
$ws.Range("N121").Value = -187505920
$ws.Range("H122").Value = 16668404
$ws.Range("I122").Value = 27778562
$ws.Range("J122").Value = 3168.75
$ws.Range("K122").Value = 250007058
$ws.Range("L122").Value = 28518.75
$ws.Range("M122").Value = -250004608
$ws.Range("N122").Value = -33418.75
$ws.Range("H123").Value = 3062.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 3062.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 9187.5
$ws.Range("N123").Value = -14087.5
$ws.Range("H124").Value = 1792.2
$ws.Range("I124").Value = 987
$ws.Range("J124").Value = 3000
$ws.Range("K124").Value = 2961
$ws.Range("L124").Value = 9000
$ws.Range("M124").Value = 1949
$ws.Range("N124").Value = -18820
$ws.Range("H125").Value = 2293.6365
$ws.Range("I125").Value = 1686
$ws.Range("J125").Value = 2800
$ws.Range("K125").Value = 5058
$ws.Range("L125").Value = 8400
$ws.Range("M125").Value = -138
$ws.Range("N125").Value = -18240
$ws.Range("H126").Value = 3667.6667
$ws.Range("I126").Value = 977.5
$ws.Range("J126").Value = 5819.8
$ws.Range("K126").Value = 2932.5
$ws.Range("L126").Value = 17459.4
$ws.Range("M126").Value = 2007.5
$ws.Range("N126").Value = -27339.4
$ws.Range("H127").Value = 1900
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1900
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 5700
$ws.Range("N127").Value = -15620
$ws.Range("H128").Value = 107457.5
$ws.Range("I128").Value = 107457.5
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 322372.5
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -317392.5
$ws.Range("H129").Value = 1166.5454
$ws.Range("I129").Value = 816.5
$ws.Range("J129").Value = 1586.6
$ws.Range("K129").Value = 2449.5
$ws.Range("L129").Value = 4759.799999999999
$ws.Range("M129").Value = 2550.5
$ws.Range("N129").Value = -14759.8
$ws.Range("H130").Value = 2433.3333
$ws.Range("I130").Value = 2250
$ws.Range("J130").Value = 3900
$ws.Range("K130").Value = 6750
$ws.Range("L130").Value = 11700
$ws.Range("M130").Value = -1730
$ws.Range("N130").Value = -21740
$ws.Range("H131").Value = 4170.552
$ws.Range("I131").Value = 5872.222
$ws.Range("J131").Value = 3404.8
$ws.Range("K131").Value = 17616.666
$ws.Range("L131").Value = 10214.4
$ws.Range("M131").Value = -12576.666
$ws.Range("N131").Value = -20294.4
$ws.Range("H132").Value = 66667520
$ws.Range("I132").Value = 142857900
$ws.Range("J132").Value = 936
$ws.Range("K132").Value = 1285721100
$ws.Range("L132").Value = 8424
$ws.Range("M132").Value = -1285718570
$ws.Range("N132").Value = -13484
$ws.Range("H133").Value = 7195.364
$ws.Range("I133").Value = 4229.8
$ws.Range("J133").Value = 9666.667
$ws.Range("K133").Value = 12689.4
$ws.Range("L133").Value = 29000.001
$ws.Range("M133").Value = -7629.400000000001
$ws.Range("N133").Value = -39120.001
$ws.Range("H134").Value = 2374.1035
$ws.Range("I134").Value = 1954
$ws.Range("J134").Value = 4999.75
$ws.Range("K134").Value = 5862
$ws.Range("L134").Value = 14999.25
$ws.Range("M134").Value = -792
$ws.Range("N134").Value = -25139.25
$ws.Range("H136").Value = 8021.25
$ws.Range("I136").Value = 2063.3333
$ws.Range("J136").Value = 10007.223
$ws.Range("K136").Value = 6189.999899999999
$ws.Range("L136").Value = 30021.669
$ws.Range("M136").Value = -1089.999899999999
$ws.Range("N136").Value = -40221.669
$ws.Range("H137").Value = 19824.213
$ws.Range("I137").Value = 2122.353
$ws.Range("J137").Value = 26663.568
$ws.Range("K137").Value = 6367.059
$ws.Range("L137").Value = 79990.704
$ws.Range("M137").Value = -1267.059
$ws.Range("N137").Value = -90190.704
$ws.Range("H138").Value = 2260.182
$ws.Range("I138").Value = 2058.7778
$ws.Range("J138").Value = 3166.5
$ws.Range("K138").Value = 6176.3334
$ws.Range("L138").Value = 9499.5
$ws.Range("M138").Value = -1036.3334
$ws.Range("N138").Value = -19779.5
$ws.Range("H139").Value = 2818.5789
$ws.Range("I139").Value = 1512
$ws.Range("J139").Value = 4270.3335
$ws.Range("K139").Value = 4536
$ws.Range("L139").Value = 12811.0005
$ws.Range("M139").Value = 604
$ws.Range("N139").Value = -23091.0005
$ws.Range("H140").Value = 1370.5938
$ws.Range("I140").Value = 719.5417
$ws.Range("J140").Value = 3323.75
$ws.Range("K140").Value = 2158.6251
$ws.Range("L140").Value = 9971.25
$ws.Range("M140").Value = 3021.3749
$ws.Range("N140").Value = -20331.25
$ws.Range("H141").Value = 501764.75
$ws.Range("I141").Value = 667353
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 2002059
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -1996879
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 626148.75
$ws.Range("I14").Value = 834698.3
$ws.Range("K14").Value = 834698.3
$ws.Range("M14").Value = -834530.3
